$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1. Drop the extra task rows (6-12), keeping only tasks 1-5.
# ------------------------------------------------------------------
$ws.Range("9:15").EntireRow.Delete() | Out-Null

# ------------------------------------------------------------------
# 2. Remove the "March" month block and weeks M:R (columns 13-18).
#    This shrinks the G1:R1 -> G1:L1 merge, and the L2:O2 -> L2 merge
#    automatically, and drops the P2:R2 merge (fully removed) along
#    with the March header cell and week-range cells M3:R3.
# ------------------------------------------------------------------
$ws.Range("M1:R1").EntireColumn.Delete() | Out-Null

# Re-assert the now single-cell merge on L2 so it is preserved as an
# explicit merged range (a size-1 merge collapses otherwise).
$ws.Range("L2").Merge() | Out-Null

# ------------------------------------------------------------------
# 3. Fix up the week-range text now that the chronogram starts a day
#    later (26/Dec instead of 25/Dec).
# ------------------------------------------------------------------
$ws.Range("F3").Value2 = "26/Dec - 31/Dec"

# ------------------------------------------------------------------
# 4. Clear the "Activity" column contents for all remaining task rows
#    (the task names are no longer populated when there is no start
#    week input), while keeping their formatting/style intact.
# ------------------------------------------------------------------
$ws.Range("C4:C8").ClearContents() | Out-Null

# ------------------------------------------------------------------
# 5. Update the Start/End Date text for the affected rows.
# ------------------------------------------------------------------
$ws.Range("D4").Value2 = "12/26"
$ws.Range("E7").Value2 = "01/21"
$ws.Range("D8").Value2 = "01/22"
$ws.Range("E8").Value2 = "01/28"

# ------------------------------------------------------------------
# 6. Fix up the diagonal "highlight" cells that mark which week(s)
#    each task spans. Row 7 now only spans a single week (I7), and
#    row 8 spans a single week at J8 instead of L8:O8.
# ------------------------------------------------------------------
$ws.Range("J7").Clear() | Out-Null
$ws.Range("K7").Clear() | Out-Null

$ws.Range("F4").Copy($ws.Range("J8")) | Out-Null
$ws.Range("L8").Clear() | Out-Null
$ws.Range("M8").Clear() | Out-Null
$ws.Range("N8").Clear() | Out-Null
$ws.Range("O8").Clear() | Out-Null
